$wb = $excel.ActiveWorkbook
$issues = $wb.Worksheets.Item("issues")
$assignees = $wb.Worksheets.Item("assignees")

# ---------------------------------------------------------------------------
# 1. Update the "Sample Transmission calculator" row body text (B10) so the
#    testing link points at the new developer test-guide page.
# ---------------------------------------------------------------------------
$issues.Range("B10").Value = "Make sure that inputs and outputs work sensibly, stress test with some bad inputs (e.g. letters in a numeric input)`nSee https://developer.mantidproject.org/Testing/General/SampleTransmissionCalculatorTestGuide.html"

# ---------------------------------------------------------------------------
# 2. Roster update: thomashampson and jhaigh0 swap groups, and two new
#    assignees (eurydice76, perenon) join jhaigh0's old group.
# ---------------------------------------------------------------------------
$groupA = "cailafinn, jhaigh0, MohamedAlmaki, adriazalvarez"
$groupB = "thomashampson, SilkeSchomann, sf1919, eurydice76, perenon"

$groupARows = @(2, 3, 6, 17, 18, 19, 20)
$groupBRows = @(5, 8, 12, 14, 16)

foreach ($r in $groupARows) {
    $issues.Range("D$r").Value = $groupA
}
foreach ($r in $groupBRows) {
    $issues.Range("D$r").Value = $groupB
}

# ---------------------------------------------------------------------------
# 3. Add the two new assignees to the "assignees" roster sheet, with the
#    same COUNTIF formula pattern used by the existing rows.
# ---------------------------------------------------------------------------
$assignees.Range("A13").Value = "eurydice76"
$assignees.Range("B13").Formula = '=COUNTIF(issues!$D$2:$D$20,"*"&A13&"*")'

$assignees.Range("A14").Value = "perenon"
$assignees.Range("B14").Formula = '=COUNTIF(issues!$D$2:$D$20,"*"&A14&"*")'

# ---------------------------------------------------------------------------
# 4. View/selection bookkeeping: "issues" becomes the active sheet/tab, with
#    B10 selected; "assignees" is left with F15 selected.
# ---------------------------------------------------------------------------
$assignees.Range("F15").Select()

$issues.Activate()
$issues.Range("B10").Select()
